$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refactor the "total" template-placeholder cells from the old
# excess-key style ("append"/"monetary"/"fn") to the new "post" function
# array style, per the commit message.
#
# Old (H10, H11, E13):
#   {{{"id": "total", "title": "Сумма", "value": "2500", "append": " руб.", "monetary": " True"}}}
# New (H10, H11, E13):
#   {{{"id": "total", "title": "Сумма", "value": "2500", "post": [{"fn": "ru_monetary_string_replace"}, {"fn": "append", "args": [" руб."]}]}}}
#
# Old (C14):
#   {{{"id": "total", "title": "Сумма", "value": "2500", "fn": "num2text"}}}
# New (C14):
#   {{{"id": "total", "title": "Сумма", "value": "2500", "post": [{"fn": "ru_monetary_as_string"}]}}}

$totalPost = '{{{"id": "total", "title": "Сумма", "value": "2500", "post": [{"fn": "ru_monetary_string_replace"}, {"fn": "append", "args": [" руб."]}]}}}'
$totalAsString = '{{{"id": "total", "title": "Сумма", "value": "2500", "post": [{"fn": "ru_monetary_as_string"}]}}}'

$ws.Cells.Item(10, 8).Value = $totalPost
$ws.Cells.Item(11, 8).Value = $totalPost
$ws.Cells.Item(14, 3).Value = $totalAsString
$ws.Cells.Item(13, 5).Value = $totalPost

# --- Touch the B9:G9 merged region (unmerge + re-merge) so it is re-appended
# to the end of the merge-cell list.
$ws.Range("B9:G9").UnMerge()
$ws.Range("B9:G9").Merge()

# --- Move the active selection from B3:H3 to the single cell G13.
$ws.Range("G13").Select()
